$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data rows (rows 2-6) that will be replaced
$ws.Range("A2:E6").ClearContents()

# Row 2
$ws.Cells.Item(2, 1).Value = "conv_7"
$ws.Cells.Item(2, 2).Value = 1689150986506
$ws.Cells.Item(2, 3).Value = 1689150989707
$ws.Cells.Item(2, 4).Value = 3
$ws.Cells.Item(2, 5).Value = "sd"

# Row 3
$ws.Cells.Item(3, 1).Value = "conv_5"
$ws.Cells.Item(3, 2).Value = 1689150991172
$ws.Cells.Item(3, 3).Value = 1689150993837
$ws.Cells.Item(3, 4).Value = 4
$ws.Cells.Item(3, 5).Value = "dxw"

# Row 4
$ws.Cells.Item(4, 1).Value = "conversation_11_07_2023__09_38_05"
$ws.Cells.Item(4, 2).Value = 1689150995253
$ws.Cells.Item(4, 3).Value = 1689150998124
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = "ded"

# Rows 5 and 6 no longer exist in the new data range; make sure they're empty
$ws.Range("A5:E6").ClearContents()

$wb.Save()
